$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New case-data rows appended to the table (rows 922-937), matching the
# target diff. Columns H, I, J, K hold numeric-looking text (e.g. "$ 0",
# "1") that must be stored as literal text, not converted to numbers, so
# those cells are pre-formatted as Text ("@") before assignment.
$data = @(
    @("21CRB01291","Hemmeter","PERMISSION REQ'D TO USE LICENSED DOCK","1501:46-12-04","MM","No Contest","Guilty","$ 0","$ 0","1","None"),
    @("21CRB01291","Hemmeter","PERMISSION REQ'D TO USE LICENSED DOCK","1501:46-12-04","MM","No Contest","Guilty","$ 0","$ 0","1","None"),
    @("21CRB01291","Hemmeter","PERMISSION REQ'D TO USE LICENSED DOCK","1501:46-12-04","MM","No Contest","Guilty","$ 0","$ 0","1","None"),
    @("21CRB01291","Hemmeter","PERMISSION REQ'D TO USE LICENSED DOCK","1501:46-12-04","MM","No Contest","Guilty","$ 0","$ 0","1","None"),
    @("21CRB01291","Hemmeter","PERMISSION REQ'D TO USE LICENSED DOCK","1501:46-12-04","MM","No Contest","Guilty","$ 0","$ 0","1","None"),
    @("21CRB01291","Hemmeter","PERMISSION REQ'D TO USE LICENSED DOCK","1501:46-12-04","MM","No Contest","Guilty","$ 0","$ 0","1","None"),
    @("21CRB01291","Hemmeter","PERMISSION REQ'D TO USE LICENSED DOCK","1501:46-12-04","MM","No Contest","Guilty","$ 0","$ 0","1","None"),
    @("21CRB01268","Hemmeter","POSSESSION DRUG PARAPHERNALIA","2925.14(C)","M4","No Contest","Guilty","$ 0","$ 0","2","None"),
    @("21CRB01268","Hemmeter","POSSESSION DRUG PARAPHERNALIA","2925.14(C)","M4","No Contest","Guilty","$ 150","$ 0","5","None"),
    @("21CRB01268","Hemmeter","POSSESSION DRUG PARAPHERNALIA","2925.14(C)","M4","No Contest","Guilty","$ 150","$ 0","5","None"),
    @("21CRB01268","Hemmeter","POSSESSION DRUG PARAPHERNALIA","2925.14(C)","M4","No Contest","Guilty","$ 150","$ 0","5","None"),
    @("21CRB01268","Hemmeter","POSSESSION DRUG PARAPHERNALIA","2925.14(C)","M4","No Contest","Guilty","$ 150","$ 0","5","None"),
    @("21CRB01268","Hemmeter","POSSESSION DRUG PARAPHERNALIA","2925.14(C)","M4","No Contest","Guilty","$ 150","$ 0","5","None"),
    @("21CRB01268","Hemmeter","POSSESSION DRUG PARAPHERNALIA","2925.14(C)","M4","No Contest","Guilty","$ 150","$ 0","5","None"),
    @("21CRB01268","Hemmeter","POSSESSION DRUG PARAPHERNALIA","2925.14(C)","M4","No Contest","Guilty","$ 150","$ 0","5","None"),
    @("21CRB01268","Hemmeter","POSSESSION DRUG PARAPHERNALIA","2925.14(C)","M4","No Contest","Guilty","$ 150","$ 0","5","None")
)

$startRow = 922
$textCols = @(8, 9, 10, 11)
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($textCols -contains $c) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $rowVals[$c - 1]
    }
}
